$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update flight number for row 2 and push its departure date forward.
$ws.Range("A2").Value = "VN0012"
$ws.Range("C2").Value = 45635.375

# Remove the now-stale sample rows (3-6) while keeping the date formatting
# on column C intact, matching the "fix flow create booking ui" cleanup.
$ws.Range("A3:I6").ClearContents()

# Leave the selection on A2, as saved in the authored workbook.
[void]$ws.Range("A2").Select()

$wb.Save()
